# "dto ok , form add/up/details , todo assignement"
#
# The backlog sheet had a small copy/content fix: the duration cell for the
# "modal restaurant details" row had a typo ("1 dat" instead of "1 day").
# Fixing it also removes the now-unreferenced shared string from the
# workbook's string table.
#
# The author's selection when the file was last saved also moved to I15.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the "1 dat" -> "1 day" typo in the real-schedule column for the
# "modal restaurant details" row (C7).
$ws.Range("C7").Value = "1 day"

# Restore the author's last selection/active cell.
$ws.Range("I15").Select()
